$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 62
$ws1.Range("F4").Value = 2011
$ws1.Range("F6").Value = 568
$ws1.Range("F8").Value = 2052
$ws1.Range("F9").Value = 10417
$ws1.Range("F11").Value = 151
$ws1.Range("F15").Value = 7276
$ws1.Range("F18").Value = 156
$ws1.Range("F19").Value = 60
$ws1.Range("F20").Value = 2435

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 62
$ws4.Range("F4").Value = 2011
$ws4.Range("F6").Value = 568
$ws4.Range("F7").Value = 19
$ws4.Range("F9").Value = 2052
$ws4.Range("F12").Value = 10417
$ws4.Range("F14").Value = 151
$ws4.Range("F18").Value = 7276
$ws4.Range("F21").Value = 156
$ws4.Range("F22").Value = 60
$ws4.Range("F23").Value = 2435
